$wb = $excel.ActiveWorkbook

# --- Sheet "Regions": remove the CA1 row (row 12) ---
$wsRegions = $wb.Worksheets.Item("Regions")
$wsRegions.Rows.Item(12).Select()
$wsRegions.Rows.Item(12).Delete()

# --- Sheet "Hierarchy": remove the CA1 row (row 10) ---
$wsHierarchy = $wb.Worksheets.Item("Hierarchy")
$wsHierarchy.Rows.Item(10).Select()
$wsHierarchy.Rows.Item(10).Delete()

# --- Sheet "Connection": remove the CA1 column (K) and CA1 row (11) ---
$wsConnection = $wb.Worksheets.Item("Connection")
$wsConnection.Columns.Item(11).Delete()
$wsConnection.Rows.Item(11).Delete()

# CA1's unique outgoing connection to FL is now represented via UAllo (Allocentric IO)
$wsConnection.Range("B10").Value = 1

# Select the sheet's used range and make "Connection" the active sheet/tab
$wsConnection.UsedRange.Select()
$wsConnection.Activate()
